$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.173.44'
$ws.Range('E2').Value = '  -1.94%  '
$ws.Range('D3').Value = '3.128.00'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  +0.04%  '
$c = $ws.Range('D5')
$c.Value = '''595.80'
$c.ClearFormats()
$ws.Range('E5').Value = '  -2.39%  '
$c = $ws.Range('D6')
$c.Value = '''136.93'
$c.ClearFormats()
$ws.Range('E6').Value = '  -4.54%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '3.123.24'
$ws.Range('E8').Value = '  -0.45%  '
$c = $ws.Range('D9')
$c.Value = '''0.520'
$c.ClearFormats()
$ws.Range('E9').Value = '  -1.81%  '
$c = $ws.Range('D10')
$c.Value = '''0.147'
$c.ClearFormats()
$ws.Range('E10').Value = '  -2.67%  '
$c = $ws.Range('D11')
$c.Value = '''5.31'
$c.ClearFormats()
$ws.Range('E11').Value = '  -0.74%  '
$c = $ws.Range('D12')
$c.Value = '''0.461'
$c.ClearFormats()
$ws.Range('E12').Value = '  -3.03%  '
$c = $ws.Range('D13')
$c.Value = '''0.0000250'
$c.ClearFormats()
$ws.Range('E13').Value = '  -2.54%  '
$c = $ws.Range('D14')
$c.Value = '''34.37'
$c.ClearFormats()
$ws.Range('E14').Value = '  -3.29%  '
$ws.Range('D15').Value = '3.637.41'
$ws.Range('E15').Value = '  -0.46%  '
$ws.Range('E16').Value = '  +1.60%  '
$ws.Range('D17').Value = '63.186.37'
$ws.Range('E17').Value = '  -1.87%  '
$ws.Range('D18').Value = '3.121.89'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').Value = '6.75'
$ws.Range('E19').Value = '  -1.68%  '
$c = $ws.Range('D20')
$c.Value = '''478.22'
$c.ClearFormats()
$ws.Range('E20').Value = '  -0.01%  '
$ws.Range('D21').Value = '14.19'
$ws.Range('E21').Value = '  -3.37%  '
$c = $ws.Range('D22')
$c.Value = '''0.701'
$c.ClearFormats()
$ws.Range('E22').Value = '  -2.84%  '
$c = $ws.Range('D23')
$c.Value = '''7.74'
$c.ClearFormats()
$ws.Range('E23').Value = '  -0.93%  '
$c = $ws.Range('D24')
$c.Value = '''87.38'
$c.ClearFormats()
$ws.Range('E24').Value = '  +2.77%  '
$ws.Range('D25').Value = '13.08'
$ws.Range('E25').Value = '  -3.80%  '
$ws.Range('E26').Value = '  +0.06%  '
$c = $ws.Range('D27')
$c.Value = '''2.72'
$c.ClearFormats()
$ws.Range('E27').Value = '  -2.18%  '
$c = $ws.Range('D28')
$c.Value = '''7.23'
$c.ClearFormats()
$ws.Range('E28').Value = '  -2.11%  '
$c = $ws.Range('D29')
$c.Value = '''8.02'
$c.ClearFormats()
$ws.Range('E29').Value = '  -6.05%  '
$c = $ws.Range('D30')
$c.Value = '''2.11'
$c.ClearFormats()
$ws.Range('E30').Value = '  +1.04%  '
$c = $ws.Range('D31')
$c.Value = '''27.17'
$c.ClearFormats()
$ws.Range('E31').Value = '  +1.80%  '
$ws.Range('E32').Value = '  -0.01%  '
$ws.Range('E33').Value = '  -7.26%  '
$c = $ws.Range('D34')
$c.Value = '''2.55'
$c.ClearFormats()
$ws.Range('E34').Value = '  -3.11%  '
$ws.Range('E35').Value = '  -2.69%  '
$c = $ws.Range('D36')
$c.Value = '''5.86'
$c.ClearFormats()
$ws.Range('E36').Value = '  -1.25%  '
$c = $ws.Range('D37')
$c.Value = '''51.98'
$c.ClearFormats()
$ws.Range('E37').Value = '  -0.95%  '
$ws.Range('D38').Value = '0.0₃0714'
$ws.Range('E38').Value = '  -3.48%  '
$c = $ws.Range('D39')
$c.Value = '''0.0392'
$c.ClearFormats()
$ws.Range('E39').Value = '  -1.14%  '
$c = $ws.Range('D40')
$c.Value = '''421.23'
$c.ClearFormats()
$ws.Range('E40').Value = '  -6.94%  '
$c = $ws.Range('D41')
$c.Value = '''0.117'
$c.ClearFormats()
$ws.Range('E41').Value = '  -0.90%  '
$c = $ws.Range('D42')
$c.Value = '''8.29'
$c.ClearFormats()
$ws.Range('E42').Value = '  -0.37%  '
$c = $ws.Range('D43')
$c.Value = '''2.70'
$c.ClearFormats()
$ws.Range('E43').Value = '  -9.48%  '
$ws.Range('D44').Value = '2.883.36'
$ws.Range('E44').Value = '  +0.92%  '
$c = $ws.Range('D45')
$c.Value = '''0.264'
$c.ClearFormats()
$ws.Range('E45').Value = '  -0.51%  '
$ws.Range('E46').Value = '  -5.38%  '
$ws.Range('E47').Value = '  -0.04%  '
$c = $ws.Range('D48')
$c.Value = '''25.78'
$c.ClearFormats()
$ws.Range('E48').Value = '  -2.32%  '
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('D50').Value = '2.29'
$ws.Range('E50').Value = '  -5.49%  '
$ws.Range('E51').Value = '  -0.62%  '
